$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Matches the new <sheetPr><outlinePr summaryBelow="0" summaryRight="0"/></sheetPr>
$ws.Outline.SummaryRow = 0
$ws.Outline.SummaryColumn = 0

# F1: replace the HYPERLINK() formula with a plain text value, keep the
# existing worksheet-level hyperlink (still pointing to the same URL).
$f1 = $ws.Range("F1")
$f1.Value = "Metadata - Single European Sky Portal"
$f1.Font.Color = 13391121
$f1.Font.Name = "Arial"

# F2: replace the HYPERLINK() mailto formula with a plain text value and
# drop the worksheet-level hyperlink that used to sit on F2 (only F2's,
# leave F1's alone).
$f2 = $ws.Range("F2")
$f2.Value = "pru-support@eurocontrol.int"

foreach ($h in @($ws.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$F$2') {
        $h.Delete()
    }
}
